$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.970.28"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "'2.356.95"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'0.687"
$ws.Range("E5").Value = "  +3.79%  "
$ws.Range("D6").Value = "'239.97"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("D7").Value = "'76.30"
$ws.Range("E7").Value = "  +4.85%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  +15.15%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'33.26"
$ws.Range("E12").Value = "  +17.57%  "
$ws.Range("D13").Value = "'7.42"
$ws.Range("E13").Value = "  +11.66%  "
$ws.Range("D14").Value = "'0.108"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "'2.707.42"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "'16.60"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "'0.912"
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").Value = "'2.358.81"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "'43.910.44"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  +5.11%  "
$ws.Range("D22").Value = "'77.46"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "'258.73"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  -2.26%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'1.86"
$ws.Range("E26").Value = "  +18.76%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.51"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").Value = "'10.78"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "'23.00"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").Value = "'174.91"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "'0.128"
$ws.Range("E32").Value = "  -3.51%  "
$ws.Range("E33").Value = "  +3.90%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0756"
$ws.Range("E34").Value = "  +6.23%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.32"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").Value = "'5.38"
$ws.Range("E36").Value = "  +4.02%  "
$ws.Range("D37").Value = "'3.73"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").Value = "'2.39"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "'6.38"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "'0.0282"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("D41").Value = "'0.212"
$ws.Range("E41").Value = "  +16.95%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.108"
$ws.Range("E42").Value = "  +11.10%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'19.36"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'9.20"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.60"
$ws.Range("E46").Value = "  +12.36%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.78"
$ws.Range("E47").Value = "  +7.98%  "
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "'100.50"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'56.75"
$ws.Range("E51").Value = "  +8.83%  "
